$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "nurse 3"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 5

$ws.Range("C6").Select()
